$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "Super, nur Tutorial ist schlecht erkennbar"
$ws.Range("L7").Value = "Weitere Animationen"
$ws.Range("E9").Value = "Super, nur Tutorial vielleicht mit blinkenden Pfeilen"
$ws.Range("L9").Value = "In Arbeit"
$ws.Range("E12").Value = "Super, Level 4 anpassen (Resetpunkt prüfen)"
$ws.Range("L12").Value = "Erledigt"
$ws.Range("L13").Value = "Erledigt"

$ws.Range("A11").Select()
